# Apply the "Updated cryptos list" data refresh to Sheet1.
# Columns: B=Coin, C=Link, D=Price, E=Volume(1h).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D/E hold numeric-looking text (e.g. "28.931.27", "  -1.69%  ") that must
# stay text, not become floating point numbers. Force Text format before
# writing, then drop back to the default style so no stray formatting is
# left behind on the cells.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value2 = "28.931.27"
$ws.Range("E2").Value2 = "  -1.69%  "
$ws.Range("D3").Value2 = "1.834.10"
$ws.Range("E3").Value2 = "  -1.93%  "
$ws.Range("D4").Value2 = "0.9995"
$ws.Range("E4").Value2 = "  -0.18%  "
$ws.Range("D5").Value2 = "245.33"
$ws.Range("E5").Value2 = "  +0.41%  "
$ws.Range("E6").Value2 = "  -1.94%  "
$ws.Range("D7").Value2 = "0.9995"
$ws.Range("E7").Value2 = "  -0.19%  "
$ws.Range("D8").Value2 = "0.07693"
$ws.Range("E8").Value2 = "  -2.62%  "
$ws.Range("D9").Value2 = "0.3052"
$ws.Range("E9").Value2 = "  -2.70%  "
$ws.Range("D10").Value2 = "23.57"
$ws.Range("E10").Value2 = "  -3.99%  "
$ws.Range("D11").Value2 = "0.07809"
$ws.Range("E11").Value2 = "  -0.68%  "
$ws.Range("D12").Value2 = "1.835.39"
$ws.Range("E12").Value2 = "  -2.14%  "
$ws.Range("D13").Value2 = "5.079"
$ws.Range("E13").Value2 = "  -2.22%  "
$ws.Range("D14").Value2 = "90.62"
$ws.Range("E14").Value2 = "  -3.46%  "
$ws.Range("D15").Value2 = "0.6806"
$ws.Range("E15").Value2 = "  -3.05%  "
$ws.Range("D16").Value2 = "6.441"
$ws.Range("E16").Value2 = "  -1.23%  "
$ws.Range("D17").Value2 = "0.000008354"
$ws.Range("E17").Value2 = "  -0.39%  "
$ws.Range("D18").Value2 = "28.943.99"
$ws.Range("E18").Value2 = "  -1.72%  "
$ws.Range("D19").Value2 = "243.36"
$ws.Range("E19").Value2 = "  -4.33%  "
$ws.Range("D20").Value2 = "2.083.41"
$ws.Range("E20").Value2 = "  -2.32%  "
$ws.Range("D21").Value2 = "12.69"
$ws.Range("E21").Value2 = "  -3.21%  "
$ws.Range("D22").Value2 = "0.9995"
$ws.Range("E22").Value2 = "  -0.15%  "
$ws.Range("D23").Value2 = "7.478"
$ws.Range("D24").Value2 = "0.9996"
$ws.Range("E24").Value2 = "  -0.20%  "
$ws.Range("D25").Value2 = "0.1469"
$ws.Range("E25").Value2 = "  -5.78%  "
$ws.Range("D26").Value2 = "162.12"
$ws.Range("E26").Value2 = "  +0.45%  "
$ws.Range("D27").Value2 = "8.802"
$ws.Range("E27").Value2 = "  -2.28%  "
$ws.Range("D28").Value2 = "18.20"
$ws.Range("E29").Value2 = "  +3.00%  "
$ws.Range("D30").Value2 = "4.219"
$ws.Range("E30").Value2 = "  -2.69%  "
$ws.Range("D31").Value2 = "4.160"
$ws.Range("E31").Value2 = "  -2.33%  "
$ws.Range("D32").Value2 = "1.180"
$ws.Range("E32").Value2 = "  -2.57%  "
$ws.Range("D33").Value2 = "0.05129"
$ws.Range("D34").Value2 = "0.7660"
$ws.Range("E34").Value2 = "  +1.94%  "
$ws.Range("D35").Value2 = "1.847"
$ws.Range("E35").Value2 = "  -2.60%  "
$ws.Range("D36").Value2 = "1.146"
$ws.Range("E36").Value2 = "  -2.55%  "
$ws.Range("D37").Value2 = "2.681"
$ws.Range("E37").Value2 = "  -1.17%  "
$ws.Range("E38").Value2 = "  -2.40%  "
$ws.Range("D39").Value2 = "1.230.18"
$ws.Range("E39").Value2 = "  -4.00%  "
$ws.Range("D40").Value2 = "2.697"
$ws.Range("E40").Value2 = "  -2.55%  "
$ws.Range("D41").Value2 = "0.9231"
$ws.Range("E41").Value2 = "  +3.05%  "
$ws.Range("D42").Value2 = "108.35"
$ws.Range("E42").Value2 = "  -0.89%  "
$ws.Range("D43").Value2 = "5.861"
$ws.Range("E43").Value2 = "  -2.77%  "
$ws.Range("D44").Value2 = "0.9989"
$ws.Range("E44").Value2 = "  -0.25%  "
$ws.Range("D45").Value2 = "9.578"
$ws.Range("E45").Value2 = "  -0.14%  "
$ws.Range("B46").Value2 = "BabyDogeCoin"
$ws.Range("C46").Value2 = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").Value2 = "0.00000000122"
$ws.Range("E46").Value2 = "  -4.53%  "
$ws.Range("B47").Value2 = "RocketPoolETH"
$ws.Range("C47").Value2 = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D47").Value2 = "1.982.26"
$ws.Range("E47").Value2 = "  -2.42%  "
$ws.Range("B48").Value2 = "Mantle"
$ws.Range("C48").Value2 = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D48").Value2 = "0.5171"
$ws.Range("E48").Value2 = "  -0.19%  "
$ws.Range("D49").Value2 = "63.99"
$ws.Range("E49").Value2 = "  -10.13%  "
$ws.Range("D50").Value2 = "1.745"
$ws.Range("E50").Value2 = "  -2.95%  "
$ws.Range("D51").Value2 = "6.934"
$ws.Range("E51").Value2 = "  -1.91%  "

$ws.Range("D2:E51").Style = "Normal"
